# ---------------------------------------------------------------------------
# Scheduled runner: refresh cached Universalis market-board price snapshots
# (currentAveragePrice / NQ / HQ) and the dependent Leve profit columns
# (LevePriceNQ/HQ, LeveProfitNQ/HQ) across the per-job Leve sheets.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 17 (Leve Item ID 38956)
$ws.Cells.Item(17, 8).Value = 1721.826   # H17
$ws.Cells.Item(17, 10).Value = 1721.826   # J17
$ws.Cells.Item(17, 12).Value = 5165.478   # L17
$ws.Cells.Item(17, 14).Value = -5501.478   # N17
# Row 18 (Leve Item ID 5471)
$ws.Cells.Item(18, 8).Value = 857.4545000000001   # H18
$ws.Cells.Item(18, 9).Value = 909   # I18
$ws.Cells.Item(18, 11).Value = 909   # K18
$ws.Cells.Item(18, 13).Value = -625   # M18
# Row 32 (Leve Item ID 5484)
$ws.Cells.Item(32, 8).Value = 3044.3333   # H32
$ws.Cells.Item(32, 9).Value = 2500.5   # I32
$ws.Cells.Item(32, 10).Value = 3199.7144   # J32
$ws.Cells.Item(32, 11).Value = 2500.5   # K32
$ws.Cells.Item(32, 12).Value = 3199.7144   # L32
$ws.Cells.Item(32, 13).Value = -2174.5   # M32
$ws.Cells.Item(32, 14).Value = -3851.7144   # N32
# Row 40 (Leve Item ID 5505)
$ws.Cells.Item(40, 8).Value = 3349.375   # H40
$ws.Cells.Item(40, 9).Value = 2975   # I40
$ws.Cells.Item(40, 10).Value = 3574   # J40
$ws.Cells.Item(40, 11).Value = 2975   # K40
$ws.Cells.Item(40, 12).Value = 3574   # L40
$ws.Cells.Item(40, 13).Value = -2800   # M40
$ws.Cells.Item(40, 14).Value = -3924   # N40
# Row 64 (Leve Item ID 5506)
$ws.Cells.Item(64, 8).Value = 7775   # H64
$ws.Cells.Item(64, 10).Value = 7999.6665   # J64
$ws.Cells.Item(64, 12).Value = 7999.6665   # L64
$ws.Cells.Item(64, 14).Value = -8495.666499999999   # N64
# Row 67 (Leve Item ID 5506)
$ws.Cells.Item(67, 8).Value = 7775   # H67
$ws.Cells.Item(67, 10).Value = 7999.6665   # J67
$ws.Cells.Item(67, 12).Value = 7999.6665   # L67
$ws.Cells.Item(67, 14).Value = -9715.666499999999   # N67
# Row 113 (Leve Item ID 27775)
$ws.Cells.Item(113, 8).Value = 2897.5652   # H113
$ws.Cells.Item(113, 9).Value = 2825.0715   # I113
$ws.Cells.Item(113, 11).Value = 2825.0715   # K113
$ws.Cells.Item(113, 13).Value = 428.9285   # M113
# Row 116 (Leve Item ID 27778)
$ws.Cells.Item(116, 8).Value = 3899.6667   # H116
$ws.Cells.Item(116, 10).Value = 3849.5   # J116
$ws.Cells.Item(116, 12).Value = 3849.5   # L116
$ws.Cells.Item(116, 14).Value = -10733.5   # N116
# Row 138 (Leve Item ID 44169)
$ws.Cells.Item(138, 8).Value = 3828.4482   # H138
$ws.Cells.Item(138, 9).Value = 3816.6667   # I138
$ws.Cells.Item(138, 10).Value = 3829.8076   # J138
$ws.Cells.Item(138, 11).Value = 11450.0001   # K138
$ws.Cells.Item(138, 12).Value = 11489.4228   # L138
$ws.Cells.Item(138, 13).Value = -6310.000100000001   # M138
$ws.Cells.Item(138, 14).Value = -21769.4228   # N138

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2 (Leve Item ID 27713)
$ws.Cells.Item(2, 8).Value = 4330230   # H2
$ws.Cells.Item(2, 9).Value = 5683190.5   # I2
$ws.Cells.Item(2, 11).Value = 5683190.5   # K2
$ws.Cells.Item(2, 13).Value = -5683077.5   # M2
# Row 74 (Leve Item ID 44000)
$ws.Cells.Item(74, 8).Value = 3935.923   # H74
$ws.Cells.Item(74, 9).Value = 1757.909   # I74
$ws.Cells.Item(74, 11).Value = 1757.909   # K74
$ws.Cells.Item(74, 13).Value = -883.9090000000001   # M74
# Row 77 (Leve Item ID 44000)
$ws.Cells.Item(77, 8).Value = 3935.923   # H77
$ws.Cells.Item(77, 9).Value = 1757.909   # I77
$ws.Cells.Item(77, 11).Value = 8789.545   # K77
$ws.Cells.Item(77, 13).Value = -4421.545   # M77
# Row 97 (Leve Item ID 19941)
$ws.Cells.Item(97, 8).Value = 1425188.4   # H97
$ws.Cells.Item(97, 9).Value = 1425188.4   # I97
$ws.Cells.Item(97, 11).Value = 1425188.4   # K97
$ws.Cells.Item(97, 13).Value = -1424692.4   # M97
# Row 116 (Leve Item ID 27713)
$ws.Cells.Item(116, 8).Value = 4330230   # H116
$ws.Cells.Item(116, 9).Value = 5683190.5   # I116
$ws.Cells.Item(116, 11).Value = 5683190.5   # K116
$ws.Cells.Item(116, 13).Value = -5680896.5   # M116
# Row 122 (Leve Item ID 36168)
$ws.Cells.Item(122, 8).Value = 1972.7333   # H122
$ws.Cells.Item(122, 9).Value = 1966   # I122
$ws.Cells.Item(122, 11).Value = 5898   # K122
$ws.Cells.Item(122, 13).Value = -3448   # M122

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3 (Leve Item ID 27713)
$ws.Cells.Item(3, 8).Value = 4330230   # H3
$ws.Cells.Item(3, 9).Value = 5683190.5   # I3
$ws.Cells.Item(3, 11).Value = 5683190.5   # K3
$ws.Cells.Item(3, 13).Value = -5683076.5   # M3
# Row 94 (Leve Item ID 19939)
$ws.Cells.Item(94, 8).Value = 2756.8572   # H94
$ws.Cells.Item(94, 9).Value = 2449.75   # I94
$ws.Cells.Item(94, 11).Value = 2449.75   # K94
$ws.Cells.Item(94, 13).Value = -1998.75   # M94
# Row 134 (Leve Item ID 43998)
$ws.Cells.Item(134, 8).Value = 3051.0322   # H134
$ws.Cells.Item(134, 9).Value = 1542.3043   # I134
$ws.Cells.Item(134, 11).Value = 4626.9129   # K134
$ws.Cells.Item(134, 13).Value = -2091.9129   # M134

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Cells.Item(31, 8).Value = 142862700   # H31
$ws.Cells.Item(31, 9).Value = 500000500   # I31
$ws.Cells.Item(31, 10).Value = 7580   # J31
$ws.Cells.Item(31, 11).Value = 500000500   # K31
$ws.Cells.Item(31, 12).Value = 7580   # L31
$ws.Cells.Item(31, 13).Value = -500000205   # M31
$ws.Cells.Item(31, 14).Value = -8170   # N31
# Row 34 (Leve Item ID 44023)
$ws.Cells.Item(34, 8).Value = 142862700   # H34
$ws.Cells.Item(34, 9).Value = 500000500   # I34
$ws.Cells.Item(34, 10).Value = 7580   # J34
$ws.Cells.Item(34, 11).Value = 500000500   # K34
$ws.Cells.Item(34, 12).Value = 7580   # L34
$ws.Cells.Item(34, 13).Value = -500000298   # M34
$ws.Cells.Item(34, 14).Value = -7984   # N34

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 41 (Leve Item ID 4700)
$ws.Cells.Item(41, 8).Value = 101.5   # H41
$ws.Cells.Item(41, 9).Value = 98.833336   # I41
$ws.Cells.Item(41, 10).Value = 109.5   # J41
$ws.Cells.Item(41, 11).Value = 296.500008   # K41
$ws.Cells.Item(41, 12).Value = 328.5   # L41
$ws.Cells.Item(41, 13).Value = 41.49999200000002   # M41
$ws.Cells.Item(41, 14).Value = -1004.5   # N41

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 69 (Leve Item ID 11891)
$ws.Cells.Item(69, 8).Value = 0   # H69
$ws.Cells.Item(69, 10).Value = 0   # J69
$ws.Cells.Item(69, 12).Value = 0   # L69
$ws.Cells.Item(69, 14).Value = $null   # N69 (cleared)
# Row 72 (Leve Item ID 11891)
$ws.Cells.Item(72, 8).Value = 0   # H72
$ws.Cells.Item(72, 10).Value = 0   # J72
$ws.Cells.Item(72, 12).Value = 0   # L72
$ws.Cells.Item(72, 14).Value = $null   # N72 (cleared)
# Row 97 (Leve Item ID 19940)
$ws.Cells.Item(97, 8).Value = 1141   # H97
$ws.Cells.Item(97, 9).Value = 1136.3   # I97
$ws.Cells.Item(97, 11).Value = 1136.3   # K97
$ws.Cells.Item(97, 13).Value = -640.3   # M97
# Row 134 (Leve Item ID 42064)
$ws.Cells.Item(134, 8).Value = 136070.58   # H134
$ws.Cells.Item(134, 10).Value = 136070.58   # J134
$ws.Cells.Item(134, 12).Value = 408211.74   # L134
$ws.Cells.Item(134, 14).Value = -413281.74   # N134

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Cells.Item(7, 8).Value = 60949.42   # H7
$ws.Cells.Item(7, 9).Value = 76582.36   # I7
$ws.Cells.Item(7, 11).Value = 76582.36   # K7
$ws.Cells.Item(7, 13).Value = -76470.36   # M7
# Row 22 (Leve Item ID 5277)
$ws.Cells.Item(22, 8).Value = 2214.9153   # H22
$ws.Cells.Item(22, 10).Value = 3967.3333   # J22
$ws.Cells.Item(22, 12).Value = 3967.3333   # L22
$ws.Cells.Item(22, 14).Value = -4557.3333   # N22
# Row 27 (Leve Item ID 5277)
$ws.Cells.Item(27, 8).Value = 2214.9153   # H27
$ws.Cells.Item(27, 10).Value = 3967.3333   # J27
$ws.Cells.Item(27, 12).Value = 3967.3333   # L27
$ws.Cells.Item(27, 14).Value = -4181.3333   # N27
# Row 46 (Leve Item ID 5282)
$ws.Cells.Item(46, 8).Value = 4157.2896   # H46
$ws.Cells.Item(46, 9).Value = 1746.6666   # I46
$ws.Cells.Item(46, 10).Value = 5729.4346   # J46
$ws.Cells.Item(46, 11).Value = 1746.6666   # K46
$ws.Cells.Item(46, 12).Value = 5729.4346   # L46
$ws.Cells.Item(46, 13).Value = -1558.6666   # M46
$ws.Cells.Item(46, 14).Value = -6105.4346   # N46
# Row 108 (Leve Item ID 25655)
$ws.Cells.Item(108, 8).Value = 35000   # H108
$ws.Cells.Item(108, 10).Value = 35000   # J108
$ws.Cells.Item(108, 12).Value = 35000   # L108
$ws.Cells.Item(108, 14).Value = -42680   # N108
# Row 126 (Leve Item ID 36249)
$ws.Cells.Item(126, 8).Value = 60949.42   # H126
$ws.Cells.Item(126, 9).Value = 76582.36   # I126
$ws.Cells.Item(126, 11).Value = 229747.08   # K126
$ws.Cells.Item(126, 13).Value = -227277.08   # M126
# Row 130 (Leve Item ID 34729)
$ws.Cells.Item(130, 8).Value = 35000   # H130
$ws.Cells.Item(130, 10).Value = 35000   # J130
$ws.Cells.Item(130, 12).Value = 35000   # L130
$ws.Cells.Item(130, 14).Value = -45040   # N130

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 34 (Leve Item ID 3349)
$ws.Cells.Item(34, 8).Value = 26333   # H34
$ws.Cells.Item(34, 9).Value = 12000   # I34
$ws.Cells.Item(34, 11).Value = 12000   # K34
$ws.Cells.Item(34, 13).Value = -11797   # M34
# Row 63 (Leve Item ID 10824)
$ws.Cells.Item(63, 8).Value = 30535.143   # H63
$ws.Cells.Item(63, 9).Value = 28226   # I63
$ws.Cells.Item(63, 10).Value = 30920   # J63
$ws.Cells.Item(63, 11).Value = 28226   # K63
$ws.Cells.Item(63, 12).Value = 30920   # L63
$ws.Cells.Item(63, 13).Value = -27602   # M63
$ws.Cells.Item(63, 14).Value = -32168   # N63
# Row 64 (Leve Item ID 11036)
$ws.Cells.Item(64, 8).Value = 49997.5   # H64
$ws.Cells.Item(64, 10).Value = 49997.5   # J64
$ws.Cells.Item(64, 12).Value = 49997.5   # L64
$ws.Cells.Item(64, 14).Value = -50493.5   # N64
# Row 66 (Leve Item ID 10824)
$ws.Cells.Item(66, 8).Value = 30535.143   # H66
$ws.Cells.Item(66, 9).Value = 28226   # I66
$ws.Cells.Item(66, 10).Value = 30920   # J66
$ws.Cells.Item(66, 11).Value = 84678   # K66
$ws.Cells.Item(66, 12).Value = 92760   # L66
$ws.Cells.Item(66, 13).Value = -81558   # M66
$ws.Cells.Item(66, 14).Value = -99000   # N66
# Row 67 (Leve Item ID 11036)
$ws.Cells.Item(67, 8).Value = 49997.5   # H67
$ws.Cells.Item(67, 10).Value = 49997.5   # J67
$ws.Cells.Item(67, 12).Value = 49997.5   # L67
$ws.Cells.Item(67, 14).Value = -51713.5   # N67
# Row 69 (Leve Item ID 10951)
$ws.Cells.Item(69, 8).Value = 34211.832   # H69
$ws.Cells.Item(69, 10).Value = 34211.832   # J69
$ws.Cells.Item(69, 12).Value = 34211.832   # L69
$ws.Cells.Item(69, 14).Value = -35709.832   # N69
# Row 70 (Leve Item ID 11979)
$ws.Cells.Item(70, 8).Value = 34000   # H70
$ws.Cells.Item(70, 9).Value = 18000   # I70
$ws.Cells.Item(70, 10).Value = 50000   # J70
$ws.Cells.Item(70, 11).Value = 18000   # K70
$ws.Cells.Item(70, 12).Value = 50000   # L70
$ws.Cells.Item(70, 13).Value = -17685   # M70
$ws.Cells.Item(70, 14).Value = -50630   # N70
# Row 72 (Leve Item ID 10951)
$ws.Cells.Item(72, 8).Value = 34211.832   # H72
$ws.Cells.Item(72, 10).Value = 34211.832   # J72
$ws.Cells.Item(72, 12).Value = 102635.496   # L72
$ws.Cells.Item(72, 14).Value = -110123.496   # N72
# Row 73 (Leve Item ID 11979)
$ws.Cells.Item(73, 8).Value = 34000   # H73
$ws.Cells.Item(73, 9).Value = 18000   # I73
$ws.Cells.Item(73, 10).Value = 50000   # J73
$ws.Cells.Item(73, 11).Value = 18000   # K73
$ws.Cells.Item(73, 12).Value = 50000   # L73
$ws.Cells.Item(73, 13).Value = -16908   # M73
$ws.Cells.Item(73, 14).Value = -52184   # N73
# Row 75 (Leve Item ID 11957)
$ws.Cells.Item(75, 8).Value = 0   # H75
$ws.Cells.Item(75, 10).Value = 0   # J75
$ws.Cells.Item(75, 12).Value = 0   # L75
$ws.Cells.Item(75, 14).Value = $null   # N75 (cleared)
# Row 78 (Leve Item ID 11957)
$ws.Cells.Item(78, 8).Value = 0   # H78
$ws.Cells.Item(78, 10).Value = 0   # J78
$ws.Cells.Item(78, 12).Value = 0   # L78
$ws.Cells.Item(78, 14).Value = $null   # N78 (cleared)
# Row 80 (Leve Item ID 10911)
$ws.Cells.Item(80, 8).Value = 0   # H80
$ws.Cells.Item(80, 9).Value = 0   # I80
$ws.Cells.Item(80, 11).Value = 0   # K80
$ws.Cells.Item(80, 13).Value = $null   # M80 (cleared)
# Row 81 (Leve Item ID 12596)
$ws.Cells.Item(81, 8).Value = 5947.143   # H81
$ws.Cells.Item(81, 9).Value = 1228.5714   # I81
$ws.Cells.Item(81, 11).Value = 2457.1428   # K81
$ws.Cells.Item(81, 13).Value = -1396.1428   # M81
# Row 83 (Leve Item ID 10911)
$ws.Cells.Item(83, 8).Value = 0   # H83
$ws.Cells.Item(83, 9).Value = 0   # I83
$ws.Cells.Item(83, 11).Value = 0   # K83
$ws.Cells.Item(83, 13).Value = $null   # M83 (cleared)
# Row 84 (Leve Item ID 12596)
$ws.Cells.Item(84, 8).Value = 5947.143   # H84
$ws.Cells.Item(84, 9).Value = 1228.5714   # I84
$ws.Cells.Item(84, 11).Value = 12285.714   # K84
$ws.Cells.Item(84, 13).Value = -6981.714   # M84
# Row 126 (Leve Item ID 36210)
$ws.Cells.Item(126, 8).Value = 1931.1111   # H126
$ws.Cells.Item(126, 9).Value = 2003.5294   # I126
$ws.Cells.Item(126, 10).Value = 700   # J126
$ws.Cells.Item(126, 11).Value = 6010.5882   # K126
$ws.Cells.Item(126, 12).Value = 2100   # L126
$ws.Cells.Item(126, 13).Value = -3540.5882   # M126
$ws.Cells.Item(126, 14).Value = -7040   # N126
# Row 132 (Leve Item ID 44029)
$ws.Cells.Item(132, 8).Value = 8097.095   # H132
$ws.Cells.Item(132, 9).Value = 6094.1   # I132
$ws.Cells.Item(132, 10).Value = 9918   # J132
$ws.Cells.Item(132, 11).Value = 18282.3   # K132
$ws.Cells.Item(132, 12).Value = 29754   # L132
$ws.Cells.Item(132, 13).Value = -15752.3   # M132
$ws.Cells.Item(132, 14).Value = -34814   # N132
